$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4COM04_GenomeAssembly")

# Update table/header column name
$ws.Range("AL1").Value = "Output [Data]"

# Add example values to the data row (row 2)
$ws.Range("B2").Value = "SAMN00000000"
$ws.Range("E2").Value = "CLC Genomics Workbench"
$ws.Range("H2").Value = "v11.0.1"
$ws.Range("K2").Value = "quality limit: 0.05"
$ws.Range("O2").Value = "EFO"
$ws.Range("P2").Value = "https://bioregistry.io/EFO:0004202"
$ws.Range("Q2").Value = "SPAdes"
$ws.Range("T2").Value = "3.15.4"
$ws.Range("Z2").Value = "10x"
$ws.Range("AC2").Value = "full"
$ws.Range("AF2").Value = "B74 RefGen_v4"
$ws.Range("AI2").Value = "txt"
